$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the date (D) and volume (J) values between row 2 and row 5
$d2 = $ws.Range("D2").Value2
$d5 = $ws.Range("D5").Value2
$j2 = $ws.Range("J2").Value2
$j5 = $ws.Range("J5").Value2

$ws.Range("D2").Value2 = $d5
$ws.Range("D5").Value2 = $d2

$ws.Range("J2").Value2 = $j5
$ws.Range("J5").Value2 = $j2
